# Regenerate save_data: update column G ("K") values with recalculated
# strikeout-based s_vals (std/mean regen) for rows 2-74 on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 2
    6  = 1
    7  = 2
    8  = 0
    9  = 1
    10 = 2
    11 = 1
    12 = 0
    13 = 3
    14 = 0
    15 = 1
    16 = 0
    17 = 0
    18 = 1
    19 = 2
    20 = 2
    21 = 1
    22 = 4
    23 = 0
    24 = 2
    25 = 2
    26 = 1
    27 = 1
    28 = 0
    29 = 0
    30 = 2
    31 = 3
    32 = 1
    33 = 1
    34 = 1
    35 = 3
    36 = 2
    37 = 0
    38 = 2
    39 = 1
    40 = 2
    41 = 1
    42 = 1
    43 = 3
    44 = 2
    45 = 1
    46 = 2
    47 = 2
    48 = 0
    49 = 0
    50 = 1
    51 = 0
    52 = 2
    53 = 0
    54 = 2
    55 = 2
    56 = 1
    57 = 1
    58 = 1
    59 = 0
    60 = 1
    61 = 3
    62 = 2
    63 = 1
    64 = 2
    65 = 2
    66 = 1
    67 = 2
    68 = 0
    69 = 0
    70 = 1
    71 = 0
    72 = 0
    73 = 2
    74 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
